$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.381.99"
$ws.Range("E2").Value = "  +4.39%  "

$ws.Range("D3").Value = "1.712.05"
$ws.Range("E3").Value = "  +1.51%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.11"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.30"
$ws.Range("E8").Value = "  +4.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.266"
$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("D12").Value = "1.956.44"

$ws.Range("D13").Value = "1.715.72"
$ws.Range("E13").Value = "  +1.79%  "

$ws.Range("E14").Value = "  +0.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.560"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.35"
$ws.Range("E16").Value = "  +0.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "253.62"
$ws.Range("E17").Value = "  +7.34%  "

$ws.Range("D18").Value = "28.315.34"
$ws.Range("E18").Value = "  +4.15%  "

$ws.Range("D19").Value = "0.0₃0747"
$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.75"
$ws.Range("E20").Value = "  -3.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("E24").Value = "  -1.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.17"
$ws.Range("E25").Value = "  +0.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.37"
$ws.Range("E26").Value = "  +0.41%  "

$ws.Range("E27").Value = "  +1.17%  "

$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("E30").Value = "  +0.99%  "

$ws.Range("E31").Value = "  +2.88%  "

$ws.Range("E32").Value = "  +0.60%  "

$ws.Range("D33").Value = "1.483.60"
$ws.Range("E33").Value = "  -3.94%  "

$ws.Range("E34").Value = "  -1.33%  "

$ws.Range("E35").Value = "  -2.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.964"
$ws.Range("E36").Value = "  +1.83%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.596"
$ws.Range("E38").Value = "  -1.47%  "

$ws.Range("E39").Value = "  +0.56%  "

$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("E41").Value = "  +0.66%  "

$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("D44").Value = "1.859.73"
$ws.Range("E44").Value = "  +1.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.26"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.802"
$ws.Range("E46").Value = "  +1.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.73"
$ws.Range("E47").Value = "  +7.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "90.01"
$ws.Range("E48").Value = "  -0.20%  "

$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("E50").Value = "  -0.62%  "

$ws.Range("E51").Value = "  -3.00%  "
